# Auto-generated: update crypto price (D) and volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.629.22"
$ws.Range("E2").Value = "  +4.35%  "

$ws.Range("D3").Value = "2.338.13"
$ws.Range("E3").Value = "  +2.34%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "548.33"
$ws.Range("E5").Value = "  +2.52%  "

$ws.Range("D6").Value = "131.96"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.95%  "

$ws.Range("D9").Value = "2.335.87"
$ws.Range("E9").Value = "  +2.24%  "

$ws.Range("E10").Value = "  +1.59%  "

$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").Value = "23.89"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").Value = "2.753.87"
$ws.Range("E15").Value = "  +2.25%  "

$ws.Range("D16").Value = "60.569.93"
$ws.Range("E16").Value = "  +4.36%  "

$ws.Range("E17").Value = "  +1.45%  "

$ws.Range("D18").Value = "2.339.48"
$ws.Range("E18").Value = "  +1.66%  "

$ws.Range("E19").Value = "  +1.36%  "

$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "316.02"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("E22").Value = "  +3.55%  "

$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").Value = "64.34"
$ws.Range("E24").Value = "  +1.96%  "

$ws.Range("E25").Value = "  +1.50%  "

$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +2.66%  "

$ws.Range("D27").Value = "7.93"
$ws.Range("E27").Value = "  -0.44%  "

$ws.Range("E28").Value = "  +7.38%  "

$ws.Range("E29").Value = "  +12.60%  "

$ws.Range("D30").Value = "173.43"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("E31").Value = "  +2.80%  "

$ws.Range("E32").Value = "  +2.42%  "

$ws.Range("E33").Value = "  +3.50%  "

$ws.Range("E34").Value = "  +11.72%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "17.96"
$ws.Range("E36").Value = "  +0.62%  "

$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("E39").Value = "  +5.33%  "

$ws.Range("D40").Value = "330.32"
$ws.Range("E40").Value = "  +14.49%  "

$ws.Range("E41").Value = "  +3.85%  "

$ws.Range("D42").Value = "38.12"
$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").Value = "140.31"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("D44").Value = "3.50"
$ws.Range("E44").Value = "  +1.76%  "

$ws.Range("D45").Value = "0.0947"
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").Value = "19.38"
$ws.Range("E46").Value = "  +7.22%  "

$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").Value = "0.0₆0224"
$ws.Range("E49").Value = "  +21.83%  "

$ws.Range("D50").Value = "0.0215"
$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("E51").Value = "  +0.74%  "
